$d = $word.ActiveDocument

# The trailing paragraph currently reads (across three runs):
#   "I would like to look for future employers" | "." | " Also, I would like to look for programming lovers. "
# The target state merges the second and third runs into one:
#   "I would like to look for future employers" | ". Also, I would like to look for programming lovers. "
# Re-finding/replacing the text that lives in the last run (leaving the first
# run untouched) causes the adjacent identically-formatted runs covered by
# the match to coalesce into a single run, which is exactly the edit the
# diff describes.
$d.Content.Find.Execute(
    "programming lovers",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "programming lovers", 2
)
